$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Section header row (192): "GetRawVotes" ---
$ws.Range("A192").Value = "GetRawVotes"
$ws.Range("A192:C192").Font.Bold = $true
$ws.Range("A192:C192").HorizontalAlignment = -4131
$ws.Range("A192:C192").Merge()

# --- Data rows (193-197) for the GetRawVotes error codes ---
# (shared-string insertion order matches the original authoring order,
#  which filled row 197's message before row 196's)
$ws.Range("A193").Value = 2101
$ws.Range("B193").Value = "CustomerId cannot be null or empty string."
$ws.Range("C193").Value = "GetRawVotes"

$ws.Range("A194").Value = 2102
$ws.Range("B194").Value = "QSetId cannot be null or empty string."
$ws.Range("C194").Value = "GetRawVotes"

$ws.Range("A195").Value = 2103
$ws.Range("B195").Value = "QSeq cannot be null or less than 1."
$ws.Range("C195").Value = "GetRawVotes"

$ws.Range("A197").Value = 2105
$ws.Range("B197").Value = "LangId Is Null Or Empty String."
$ws.Range("C197").Value = "GetRawVotes"

$ws.Range("A196").Value = 2104
$ws.Range("B196").Value = "Begin Date and End Date cannot be null."
$ws.Range("C196").Value = "GetRawVotes"

# --- Update view selection to match the end of the new data ---
$ws.Range("C197").Select()
